# Update Skeena chinook escapement row (C12) from "Total Esc (2021)" to "GSI esc (2022)",
# and fill in the previously-blank H17 (un-merging H16:H17) with matching values/formats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C12: TCCHINOOK Table B3 escapement note text -------------------------
$newText = "TCCHINOOK Table B3: NBC Escapement  (Area 4: Skeena R. GSI esc) (2022)"
$cell = $ws.Range("C12")
$cell.Value = $newText

# Re-apply the rich-text emphasis that existed on the original note:
#  - "TCCHINOOK Table B3" prefix is bold
#  - the trailing "(2022)" is bold + blue, using the smaller "(Body)" font
$run1Len = ("TCCHINOOK Table B3").Length
$run1 = $cell.Characters(1, $run1Len)
$run1.Font.Bold = $true
$run1.Font.Size = 12
$run1.Font.Name = "Aptos Narrow"

$run3Text = "(2022)"
$run3Start = $newText.Length - $run3Text.Length + 1
$run3 = $cell.Characters($run3Start, $run3Text.Length)
$run3.Font.Bold = $true
$run3.Font.Size = 12
$run3.Font.Name = "Aptos Narrow (Body)"
$run3.Font.Color = 16724484

# --- H16 / H17: split the merged "No data available" cell -----------------
# Previously H16:H17 was merged as a single centred "No data available" cell.
# Now H16 gets its own "Expanded from spawner surveys (2022)" note (matching
# the D/E/F/G columns on the same row) and H17 keeps "No data available" as
# its own standalone cell matching the D17:G17 styling.
$ws.Range("H16:H17").UnMerge()

$ws.Range("D16").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H16").Value = "Expanded from spawner surveys (2022)"

$ws.Range("D17").Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("H17").Value = "No data available"

$excel.CutCopyMode = $false

# --- view state: last edit focus was on C12 --------------------------------
$ws.Range("C12").Select()
